$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capacitive buton (Touch Buton) row - product code & related attributes changed
# from the old "Touch" SKU to the new "One Shot" capacitive SKU.
$ws.Range("B39").Value = "BT-B18-OS-TCH-HLD-P3B2-01"
$ws.Range("D39").Value = "One Shot"
$ws.Range("F39").Value = "Tek/Çift işlevli Led"
$ws.Range("H39").Value = "Kablolu/Vidalı"
$ws.Range("I39").Value = "Buzzerlı/Buzzersız"

# Update the displayed hyperlink text on K39 to match the new product code,
# keeping the same hyperlink target/relationship.
$ws.Range("K39").Value = "https://github.com/btk42/BT-B18-OS-TCH-HLD-P3B2-01"

# Move the current selection to reflect where the editor ended up working.
$ws.Range("B39").Select()
